$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Row 107: update Jira-id (B) and Description (C), grow the row height.
# ---------------------------------------------------------------------------
$ws.Cells.Item(107, 2).Value = "OOPQA-1226|PQA-1227"
$ws.Cells.Item(107, 3).Value = "Verify that following options get displayed in SORT BY drop down in POSTS search results page: `na)Relevance `nb)Create Date(Newest) `nc)Create Date(Oldest)`nVerify that search results are sorted by CREATE DATE(NEWEST) by default in POSTS search results page"
$ws.Rows.Item(107).RowHeight = 75

# ---------------------------------------------------------------------------
# 2) Two brand-new test cases appended as rows 117 and 118.
#    Formats are cloned from existing cells (via copy / paste-special-formats)
#    so the shared style table is reused instead of growing needlessly.
# ---------------------------------------------------------------------------

# --- Row 117 -------------------------------------------------------------
$ws.Cells.Item(117, 1).Value = "TestCase_B116"
$ws.Cells.Item(108, 1).Copy()
$ws.Cells.Item(117, 1).PasteSpecial(-4122)

$ws.Cells.Item(117, 2).Value = "OPQA-1228"
$ws.Cells.Item(107, 2).Copy()
$ws.Cells.Item(117, 2).PasteSpecial(-4122)

$ws.Cells.Item(117, 3).Value = "Verify that more search results get displayed when user scrolls down in POSTS search results page"
$ws.Cells.Item(107, 2).Copy()
$ws.Cells.Item(117, 3).PasteSpecial(-4122)

$ws.Cells.Item(117, 4).Value = "Y"
$ws.Cells.Item(116, 4).Copy()
$ws.Cells.Item(117, 4).PasteSpecial(-4122)

$ws.Cells.Item(107, 5).Copy()
$ws.Cells.Item(117, 5).PasteSpecial(-4122)

# --- Row 118 -------------------------------------------------------------
$ws.Cells.Item(118, 1).Value = "TestCase_B117"
$ws.Cells.Item(108, 1).Copy()
$ws.Cells.Item(118, 1).PasteSpecial(-4122)

$ws.Cells.Item(118, 2).Value = "OPQA-1229|OPQA-1230"
$ws.Cells.Item(107, 2).Copy()
$ws.Cells.Item(118, 2).PasteSpecial(-4122)

$ws.Cells.Item(118, 3).Value = "Verify that sorting is retained when user navigates back to POSTS search results page from record view page`nVerify that search drop down content type is retained when user navigates back to POSTS search results page from record view page"
$ws.Cells.Item(107, 2).Copy()
$ws.Cells.Item(118, 3).PasteSpecial(-4122)
$ws.Cells.Item(118, 3).WrapText = $true

$ws.Cells.Item(118, 4).Value = "Y"
$ws.Cells.Item(107, 4).Copy()
$ws.Cells.Item(118, 4).PasteSpecial(-4122)

$ws.Cells.Item(107, 5).Copy()
$ws.Cells.Item(118, 5).PasteSpecial(-4122)

$ws.Rows.Item(118).RowHeight = 30

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) Column B is a little wider now that it holds multi-id Jira strings.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 22.59

# ---------------------------------------------------------------------------
# 4) Refresh the view: scrolled a row earlier, selection parked on D125.
# ---------------------------------------------------------------------------
$win = $excel.Windows.Item(1)
$win.ScrollRow = 107
$win.ScrollColumn = 1
$ws.Range("D125").Select()
